# Timing issue fix - keywords, updated tc1,2 in ubc01
#
# The "Cases" tab query (row 2 / cell B2) previously returned an extra
# `Cohort` column (coalesce(co.cohort_description, '') AS `Cohort`).
# That trailing column is removed here, which also shrinks the amount of
# wrapped text needed for the row, so its row height shrinks accordingly.
# The workbook view/selection is updated to reflect the new focused cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "CasesTab" query text (B2) - same query as before, minus the trailing
# `Cohort` column that used to be appended to the RETURN clause.
$newCasesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
 MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
	WHERE s.clinical_study_designation IN ['UBC01'] and diag.stage_of_disease in [ 'T3N0M0', 'T3N0M1', 'T3N1M0']  OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

# Update the Cases query cell (B4/"FilesTab" row and B3/"SamplesTab" row are
# untouched - only the Cases query content changes).
$ws.Range("B2").Value = $newCasesQuery

# The shorter text needs less wrapped height - row 2 shrinks from 304.5 to 290
# (matching row 3 / row 4, which already used that height).
$ws.Rows(2).RowHeight = 290

# Focus/select the edited cell (was C4:E4 on row 4, now B2 on row 2).
$ws.Range("B2").Select()
